$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: clear the "Columna Ocupada Al Pedo" header in D1 (keep style)
# ---------------------------------------------------------------------------
$ws.Range("D1").ClearContents()

# ---------------------------------------------------------------------------
# 2. Update phone numbers for the "valid message" customers (Pablo, Pedro with
#    valid message, Jose valid message, Joao message, Sicrano) and apply the
#    new look (Calibri 10 + thin top/bottom border) to those cells.
# ---------------------------------------------------------------------------
$phoneCells = @("B2", "B3", "B6", "B12", "B18")
foreach ($addr in $phoneCells) {
    $rng = $ws.Range($addr)
    $rng.Value = "55(71) 9101-3035"
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 10
    $rng.VerticalAlignment = -4160
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(8).Color = 9406680
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(9).Weight = 2
    $rng.Borders.Item(9).Color = 9406680
}

# ---------------------------------------------------------------------------
# 3. Move the "Garbage" helper cells: drop the old ones (E3, G6) and place new
#    ones at D2 and E5.
# ---------------------------------------------------------------------------
$ws.Range("E3").ClearContents()
$ws.Range("G6").ClearContents()

$ws.Range("D2").Value = "Garbage"
$ws.Range("E5").Value = "Garbage"

# ---------------------------------------------------------------------------
# 4. New blank "spacer" cells in column D with the Consolas / vertical-center
#    look (leftover formatting from the edited rows).
# ---------------------------------------------------------------------------
$blankCells = @("D3", "D10", "D13", "D15", "D17")
foreach ($addr in $blankCells) {
    $rng = $ws.Range($addr)
    $rng.ClearContents()
    $rng.Font.Name = "Consolas"
    $rng.Font.Size = 11
    $rng.Font.Color = 7901646
    $rng.VerticalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 5. Column D formatting + dimension/selection bookkeeping.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 26.7109375
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

$ws.Range("E6").Select()
